$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, pushing existing rows 53:65 down to 54:66
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new weekly record
$ws.Cells.Item(53, 1).Value = 4
$ws.Cells.Item(53, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(53, 3).Value = "Los Lagos"
$ws.Cells.Item(53, 4).Value = 44642
$ws.Cells.Item(53, 4).NumberFormat = $ws.Cells.Item(54, 4).NumberFormat
$ws.Cells.Item(53, 5).Value = 10
$ws.Cells.Item(53, 6).Value = 100112031
$ws.Cells.Item(53, 7).Value = "Poroto verde"
$ws.Cells.Item(53, 8).Value = "Magnum"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 60
$ws.Cells.Item(53, 11).Value = 36000
$ws.Cells.Item(53, 12).Value = 36000
$ws.Cells.Item(53, 13).Value = 36000
$ws.Cells.Item(53, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(53, 15).Value = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value = 1440
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"
